$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2
$ws.Range("T2").Value = 491246

# Move the active selection from T3 to T2
$ws.Range("T2").Select()
